$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 405
$ws.Range("I33").Value = 101.72
$ws.Range("J33").Value = 2300.5
$ws.Range("K33").Value = 101.72
$ws.Range("L33").Value = 2300.5
$ws.Range("M33").Value = 127.28
$ws.Range("N33").Value = -2758.5
$ws.Range("H98").Value = 78808.2
$ws.Range("I98").Value = 913.125
$ws.Range("J98").Value = 390388.5
$ws.Range("K98").Value = 913.125
$ws.Range("L98").Value = 390388.5
$ws.Range("M98").Value = 584.875
$ws.Range("N98").Value = -393384.5
$ws.Range("H103").Value = 774.375
$ws.Range("I103").Value = 600
$ws.Range("J103").Value = 948.75
$ws.Range("K103").Value = 1800
$ws.Range("L103").Value = 2846.25
$ws.Range("M103").Value = -1214
$ws.Range("N103").Value = -4018.25
$ws.Range("H122").Value = 78808.2
$ws.Range("I122").Value = 913.125
$ws.Range("J122").Value = 390388.5
$ws.Range("K122").Value = 2739.375
$ws.Range("L122").Value = 1171165.5
$ws.Range("M122").Value = -289.375
$ws.Range("N122").Value = -1176065.5
$ws.Range("H132").Value = 20615.4
$ws.Range("I132").Value = 2774.5366
$ws.Range("J132").Value = 101890.445
$ws.Range("K132").Value = 8323.6098
$ws.Range("L132").Value = 305671.335
$ws.Range("M132").Value = -5793.6098
$ws.Range("N132").Value = -310731.335
$ws.Range("H138").Value = 2594.0952
$ws.Range("I138").Value = 2052.889
$ws.Range("K138").Value = 6158.667
$ws.Range("M138").Value = -1018.667
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3125
$ws.Range("I2").Value = 2533.3333
$ws.Range("J2").Value = 4900
$ws.Range("K2").Value = 2533.3333
$ws.Range("L2").Value = 4900
$ws.Range("M2").Value = -2420.3333
$ws.Range("N2").Value = -5126
$ws.Range("H61").Value = 2215.1794
$ws.Range("I61").Value = 1678.3334
$ws.Range("J61").Value = 3074.1333
$ws.Range("K61").Value = 1678.3334
$ws.Range("L61").Value = 3074.1333
$ws.Range("M61").Value = -1466.3334
$ws.Range("N61").Value = -3498.1333
$ws.Range("H116").Value = 3125
$ws.Range("I116").Value = 2533.3333
$ws.Range("J116").Value = 4900
$ws.Range("K116").Value = 2533.3333
$ws.Range("L116").Value = 4900
$ws.Range("M116").Value = -239.3332999999998
$ws.Range("N116").Value = -9488
$ws.Range("H132").Value = 15628280
$ws.Range("I132").Value = 35716676
$ws.Range("J132").Value = 3972.5557
$ws.Range("K132").Value = 107150028
$ws.Range("L132").Value = 11917.6671
$ws.Range("M132").Value = -107147498
$ws.Range("N132").Value = -16977.6671
$ws.Range("H136").Value = 2215.1794
$ws.Range("I136").Value = 1678.3334
$ws.Range("J136").Value = 3074.1333
$ws.Range("K136").Value = 5035.0002
$ws.Range("L136").Value = 9222.3999
$ws.Range("M136").Value = -2485.0002
$ws.Range("N136").Value = -14322.3999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3125
$ws.Range("I3").Value = 2533.3333
$ws.Range("J3").Value = 4900
$ws.Range("K3").Value = 2533.3333
$ws.Range("L3").Value = 4900
$ws.Range("M3").Value = -2419.3333
$ws.Range("N3").Value = -5128
$ws.Range("H94").Value = 1397.0769
$ws.Range("I94").Value = 1397.4546
$ws.Range("J94").Value = 1395
$ws.Range("K94").Value = 1397.4546
$ws.Range("L94").Value = 1395
$ws.Range("M94").Value = -946.4546
$ws.Range("N94").Value = -2297
$ws.Range("H105").Value = 5461.375
$ws.Range("I105").Value = 4010
$ws.Range("J105").Value = 5668.7144
$ws.Range("K105").Value = 4010
$ws.Range("L105").Value = 5668.7144
$ws.Range("M105").Value = -2263
$ws.Range("N105").Value = -9162.714400000001
$ws.Range("H134").Value = 2436.9443
$ws.Range("I134").Value = 1928.3077
$ws.Range("J134").Value = 3759.4
$ws.Range("K134").Value = 5784.9231
$ws.Range("L134").Value = 11278.2
$ws.Range("M134").Value = -3249.9231
$ws.Range("N134").Value = -16348.2
$ws.Range("H138").Value = 15800
$ws.Range("J138").Value = 15800
$ws.Range("L138").Value = 15800
$ws.Range("N138").Value = -26080
$ws.Range("H140").Value = 36799.668
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 36799.668
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 36799.668
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -47159.668
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 142890580
$ws.Range("I23").Value = 250006000
$ws.Range("J23").Value = 70010
$ws.Range("K23").Value = 250006000
$ws.Range("L23").Value = 70010
$ws.Range("M23").Value = -250005760
$ws.Range("N23").Value = -70490
$ws.Range("H27").Value = 142890580
$ws.Range("I27").Value = 250006000
$ws.Range("J27").Value = 70010
$ws.Range("K27").Value = 250006000
$ws.Range("L27").Value = 70010
$ws.Range("M27").Value = -250005808
$ws.Range("N27").Value = -70394
$ws.Range("H51").Value = 83365280
$ws.Range("J51").Value = 38331.8
$ws.Range("L51").Value = 38331.8
$ws.Range("N51").Value = -39803.8
$ws.Range("H58").Value = 12822005
$ws.Range("I58").Value = 860.53845
$ws.Range("J58").Value = 38464292
$ws.Range("K58").Value = 860.53845
$ws.Range("L58").Value = 38464292
$ws.Range("M58").Value = -657.53845
$ws.Range("N58").Value = -38464698
$ws.Range("H59").Value = 44978.75
$ws.Range("J59").Value = 44978.75
$ws.Range("L59").Value = 44978.75
$ws.Range("N59").Value = -47268.75
$ws.Range("H60").Value = 48990
$ws.Range("J60").Value = 48990
$ws.Range("L60").Value = 48990
$ws.Range("N60").Value = -50012
$ws.Range("H61").Value = 83365280
$ws.Range("J61").Value = 38331.8
$ws.Range("L61").Value = 38331.8
$ws.Range("N61").Value = -39027.8
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 15500
$ws.Range("J74").Value = 15500
$ws.Range("L74").Value = 15500
$ws.Range("N74").Value = -17248
$ws.Range("H77").Value = 15500
$ws.Range("J77").Value = 15500
$ws.Range("L77").Value = 46500
$ws.Range("N77").Value = -55236
$ws.Range("H99").Value = 1763.25
$ws.Range("I99").Value = 1480.875
$ws.Range("J99").Value = 2328
$ws.Range("K99").Value = 1480.875
$ws.Range("L99").Value = 2328
$ws.Range("M99").Value = 17.125
$ws.Range("N99").Value = -5324
$ws.Range("H126").Value = 1763.25
$ws.Range("I126").Value = 1480.875
$ws.Range("J126").Value = 2328
$ws.Range("K126").Value = 4442.625
$ws.Range("L126").Value = 6984
$ws.Range("M126").Value = -1972.625
$ws.Range("N126").Value = -11924
$ws.Range("H136").Value = 12822005
$ws.Range("I136").Value = 860.53845
$ws.Range("J136").Value = 38464292
$ws.Range("K136").Value = 2581.61535
$ws.Range("L136").Value = 115392876
$ws.Range("M136").Value = -31.61535000000003
$ws.Range("N136").Value = -115397976
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 27599
$ws.Range("I18").Value = 30328.9
$ws.Range("K18").Value = 90986.70000000001
$ws.Range("M18").Value = -90817.70000000001
$ws.Range("H86").Value = 1225.75
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1225.75
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3677.25
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -6049.25
$ws.Range("H89").Value = 1225.75
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1225.75
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 11031.75
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -22887.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 68621
$ws.Range("I126").Value = 133504
$ws.Range("J126").Value = 3738
$ws.Range("K126").Value = 400512
$ws.Range("L126").Value = 11214
$ws.Range("M126").Value = -398042
$ws.Range("N126").Value = -16154
$ws.Range("H132").Value = 2928.0625
$ws.Range("I132").Value = 2150.5557
$ws.Range("J132").Value = 3927.7144
$ws.Range("K132").Value = 6451.6671
$ws.Range("L132").Value = 11783.1432
$ws.Range("M132").Value = -3921.6671
$ws.Range("N132").Value = -16843.1432
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 13484.5
$ws.Range("I34").Value = 5000
$ws.Range("J34").Value = 21969
$ws.Range("K34").Value = 5000
$ws.Range("L34").Value = 21969
$ws.Range("M34").Value = -4828
$ws.Range("N34").Value = -22313
$ws.Range("H40").Value = 7357.5713
$ws.Range("I40").Value = 7625.25
$ws.Range("J40").Value = 7000.6665
$ws.Range("K40").Value = 7625.25
$ws.Range("L40").Value = 7000.6665
$ws.Range("M40").Value = -7489.25
$ws.Range("N40").Value = -7272.6665
$ws.Range("H55").Value = 443.17648
$ws.Range("J55").Value = 404.25
$ws.Range("L55").Value = 404.25
$ws.Range("N55").Value = -750.25
$ws.Range("H132").Value = 3984.25
$ws.Range("I132").Value = 3163.2856
$ws.Range("J132").Value = 4805.2144
$ws.Range("K132").Value = 9489.856800000001
$ws.Range("L132").Value = 14415.6432
$ws.Range("M132").Value = -6959.856800000001
$ws.Range("N132").Value = -19475.6432
$ws.Range("H134").Value = 61896.332
$ws.Range("J134").Value = 61896.332
$ws.Range("L134").Value = 61896.332
$ws.Range("N134").Value = -72036.33199999999
$ws.Range("H138").Value = 58449.5
$ws.Range("J138").Value = 58449.5
$ws.Range("L138").Value = 58449.5
$ws.Range("N138").Value = -68729.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2564.348
$ws.Range("I132").Value = 2058.5
$ws.Range("K132").Value = 6175.5
$ws.Range("M132").Value = -3645.5
$ws.Range("H136").Value = 200929.84
$ws.Range("I136").Value = 263878.66
$ws.Range("J136").Value = 1591.9166
$ws.Range("K136").Value = 791635.98
$ws.Range("L136").Value = 4775.7498
$ws.Range("M136").Value = -789085.98
$ws.Range("N136").Value = -9875.7498

Write-Host "Applied all changes"